$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Header + data values
$ws.Range("A1").Value = "p_Quantity"
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4

# Border + unlocked protection for the parameter cells A2:A4
foreach ($addr in @("A2", "A3", "A4")) {
    $rng = $ws.Range($addr)
    $rng.Borders.Color = 0
    $rng.Borders.Weight = 2
    $rng.Borders.LineStyle = 1
    $rng.Borders.Item(7).LineStyle = -4142
    $rng.Locked = $false
}

# Column A should auto-size to fit "p_Quantity"
$ws.Columns.Item(1).EntireColumn.AutoFit()

# Selection lands on A2 after the edit
$ws.Range("A2").Select()
